$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy the existing header
# formatting (bold, centered, bordered) from H1 so the new columns match
# the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF)
$data = @(
    @(2, 8, 8),
    @(3, 9, 9),
    @(4, 3, 4),
    @(5, 9, 9),
    @(6, 7, 7),
    @(7, 9, 9),
    @(8, 7, 8),
    @(9, 8, 8),
    @(10, 9, 9),
    @(11, 5, 5),
    @(12, 9, 9),
    @(13, 5, 5),
    @(14, 9, 9),
    @(15, 4, 5),
    @(16, 4, 4),
    @(17, 5, 6),
    @(18, 8, 8),
    @(19, 6, 6),
    @(20, 7, 7),
    @(21, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
